$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dataRange = $ws.Range("A6:D7")
$dataRange.NumberFormat = "@"

$ws.Range("A6").Value = "2025-09-13"
$ws.Range("B6").Value = "AAA"
$ws.Range("C6").Value = "44CDX012"
$ws.Range("D6").Value = "MAMA Nagar"

$ws.Range("A7").Value = "2025-08-13"
$ws.Range("B7").Value = "JJJ"
$ws.Range("C7").Value = "456CDX0176"
$ws.Range("D7").Value = "AMMA Nagar"
